$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "324.67", "1.002").
# Force it to remain text so the updated values keep their exact display form
# (matching the original inline-string / text cell representation).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.910.89"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "1.879.22"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "324.67"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "0.4606"
$ws.Range("E7").Value = "  -1.00%  "
$ws.Range("D8").Value = "0.3868"
$ws.Range("E8").Value = "  -1.46%  "
$ws.Range("D9").Value = "0.07847"
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("D10").Value = "0.9849"
$ws.Range("E10").Value = "  -2.75%  "
$ws.Range("D11").Value = "21.78"
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").Value = "1.880.01"
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("D13").Value = "6.988"
$ws.Range("E13").Value = "  -2.13%  "
$ws.Range("D14").Value = "5.642"
$ws.Range("E14").Value = "  -2.31%  "
$ws.Range("D15").Value = "0.06978"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").Value = "87.98"
$ws.Range("E16").Value = "  -1.56%  "
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "0.000009953"
$ws.Range("E18").Value = "  -1.95%  "
$ws.Range("D19").Value = "16.89"
$ws.Range("E19").Value = "  -2.19%  "
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "28.913.82"
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("E22").Value = "  -2.01%  "
$ws.Range("E23").Value = "  -1.61%  "
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("D25").Value = "156.32"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("D26").Value = "19.33"
$ws.Range("E26").Value = "  -2.67%  "
$ws.Range("D27").Value = "5.975"
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("D28").Value = "117.53"
$ws.Range("E28").Value = "  -2.39%  "
$ws.Range("D29").Value = "1.910"
$ws.Range("E29").Value = "  -4.06%  "
$ws.Range("D30").Value = "0.09344"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").Value = "0.9014"
$ws.Range("E31").Value = "  -4.30%  "
$ws.Range("D32").Value = "5.262"
$ws.Range("E32").Value = "  -1.96%  "
$ws.Range("E33").Value = "  -2.34%  "
$ws.Range("D34").Value = "3.253"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").Value = "1.175"
$ws.Range("E35").Value = "  -0.43%  "
$ws.Range("D36").Value = "0.05746"
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("E37").Value = "  -1.33%  "
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("D39").Value = "7.635"
$ws.Range("E39").Value = "  -6.09%  "
$ws.Range("D40").Value = "0.5644"
$ws.Range("E40").Value = "  -3.45%  "
$ws.Range("E41").Value = "  -2.89%  "
$ws.Range("E42").Value = "  -2.74%  "
$ws.Range("D43").Value = "2.256"
$ws.Range("E43").Value = "  -0.98%  "
$ws.Range("D44").Value = "11.87"
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("D45").Value = "0.5336"
$ws.Range("E45").Value = "  -2.48%  "
$ws.Range("D46").Value = "0.07044"
$ws.Range("E46").Value = "  -2.41%  "
$ws.Range("D47").Value = "1.839"
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("D48").Value = "2.542"
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("D49").Value = "112.61"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("D50").Value = "1.060"
$ws.Range("E50").Value = "  -5.40%  "
$ws.Range("D51").Value = "70.72"
$ws.Range("E51").Value = "  -0.86%  "
